$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source feed re-fetched odds for several fixtures; on the provider's
# side the rows for a couple of paired matches (same matchday) came back in
# swapped order, and a handful of already-scraped upcoming fixtures got
# refreshed odds. Re-apply both kinds of change here.

function Swap-Rows {
    param($sheet, [int]$rowA, [int]$rowB, [int]$firstCol, [int]$lastCol)
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cellA = $sheet.Cells.Item($rowA, $c)
        $cellB = $sheet.Cells.Item($rowB, $c)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

# Columns B (2) through AC (29) hold the match data; column A (the running
# id) stays put for each physical row.
$firstCol = 2
$lastCol = 29

Swap-Rows $ws 188 189 $firstCol $lastCol
Swap-Rows $ws 237 238 $firstCol $lastCol
Swap-Rows $ws 251 252 $firstCol $lastCol
Swap-Rows $ws 264 265 $firstCol $lastCol
Swap-Rows $ws 282 283 $firstCol $lastCol

# Refreshed odds for a batch of upcoming fixtures (rows 288-296).
$ws.Cells.Item(288, 14).Value2 = 2.625   # N288
$ws.Cells.Item(288, 16).Value2 = 2.7     # P288
$ws.Cells.Item(288, 18).Value2 = 1.85    # R288
$ws.Cells.Item(288, 19).Value2 = 2       # S288
$ws.Cells.Item(288, 21).Value2 = 1.875   # U288
$ws.Cells.Item(288, 22).Value2 = 1.975   # V288

$ws.Cells.Item(289, 15).Value2 = 3.6     # O289
$ws.Cells.Item(289, 16).Value2 = 4       # P289

$ws.Cells.Item(290, 16).Value2 = 5       # P290
$ws.Cells.Item(290, 18).Value2 = 1.95    # R290
$ws.Cells.Item(290, 19).Value2 = 1.9     # S290
$ws.Cells.Item(290, 21).Value2 = 1.825   # U290
$ws.Cells.Item(290, 22).Value2 = 2.025   # V290

$ws.Cells.Item(291, 15).Value2 = 3.6     # O291
$ws.Cells.Item(291, 16).Value2 = 4.333   # P291

$ws.Cells.Item(292, 16).Value2 = 6.5     # P292
$ws.Cells.Item(292, 18).Value2 = 1.875   # R292
$ws.Cells.Item(292, 19).Value2 = 1.975   # S292

$ws.Cells.Item(293, 14).Value2 = 4.2     # N293
$ws.Cells.Item(293, 18).Value2 = 1.925   # R293
$ws.Cells.Item(293, 19).Value2 = 1.925   # S293
$ws.Cells.Item(293, 21).Value2 = 1.975   # U293
$ws.Cells.Item(293, 22).Value2 = 1.875   # V293

$ws.Cells.Item(294, 15).Value2 = 3.3     # O294
$ws.Cells.Item(294, 16).Value2 = 3.4     # P294
$ws.Cells.Item(294, 18).Value2 = 1.875   # R294
$ws.Cells.Item(294, 19).Value2 = 1.975   # S294
$ws.Cells.Item(294, 21).Value2 = 1.875   # U294
$ws.Cells.Item(294, 22).Value2 = 1.975   # V294

$ws.Cells.Item(295, 14).Value2 = 1.444   # N295
$ws.Cells.Item(295, 15).Value2 = 4.5     # O295
$ws.Cells.Item(295, 16).Value2 = 7.5     # P295
$ws.Cells.Item(295, 18).Value2 = 1.975   # R295
$ws.Cells.Item(295, 19).Value2 = 1.875   # S295

$ws.Cells.Item(296, 16).Value2 = 4.333   # P296
$ws.Cells.Item(296, 17).Value2 = -0.75   # Q296
$ws.Cells.Item(296, 18).Value2 = 2.05    # R296
$ws.Cells.Item(296, 19).Value2 = 1.8     # S296
